$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 24; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 29; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 39; I = '%'; J = 'Uninterpretable' },
    @{ Row = 62; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 66; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 68; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 74; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 81; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 91; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 95; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 112; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 116; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 117; I = 'ba'; J = 'Appreciation' },
    @{ Row = 126; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 128; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 129; I = 'ba'; J = 'Appreciation' },
    @{ Row = 135; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 138; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 147; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 153; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 161; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 167; I = 'ba'; J = 'Appreciation' },
    @{ Row = 168; I = 'ba'; J = 'Appreciation' },
    @{ Row = 171; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 181; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 191; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 195; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 200; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 208; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 216; I = '%'; J = 'Uninterpretable' },
    @{ Row = 231; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 240; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 246; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 251; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 253; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 257; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 267; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 276; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 277; I = '%'; J = 'Uninterpretable' },
    @{ Row = 278; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 282; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 288; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 296; I = '%'; J = 'Uninterpretable' },
    @{ Row = 302; I = '%'; J = 'Uninterpretable' },
    @{ Row = 303; I = '%'; J = 'Uninterpretable' },
    @{ Row = 311; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 316; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 326; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 331; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 349; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 356; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 357; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 360; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 371; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 373; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 379; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 381; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 387; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 388; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 390; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 405; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 429; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"
